$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.391.69"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3
$ws.Range("D3").Value = "3.520.77"
$ws.Range("E3").Value = "  -0.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'612.41"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6
$ws.Range("D6").Value = "'151.34"
$ws.Range("E6").Value = "  -1.52%  "

# Row 7
$ws.Range("D7").Value = "3.518.57"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").Value = "'0.478"
$ws.Range("E9").Value = "  -1.51%  "

# Row 10
$ws.Range("E10").Value = "  -1.03%  "

# Row 11
$ws.Range("D11").Value = "'7.10"
$ws.Range("E11").Value = "  +2.74%  "

# Row 12
$ws.Range("E12").Value = "  -1.32%  "

# Row 13
$ws.Range("E13").Value = "  -1.47%  "

# Row 14
$ws.Range("D14").Value = "4.116.25"

# Row 15
$ws.Range("D15").Value = "'32.01"
$ws.Range("E15").Value = "  -0.37%  "

# Row 16
$ws.Range("D16").Value = "3.518.17"
$ws.Range("E16").Value = "  -0.80%  "

# Row 17
$ws.Range("D17").Value = "67.366.34"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("D19").Value = "'6.39"
$ws.Range("E19").Value = "  +0.40%  "

# Row 20
$ws.Range("E20").Value = "  -1.73%  "

# Row 21
$ws.Range("D21").Value = "'444.63"
$ws.Range("E21").Value = "  -2.00%  "

# Row 22
$ws.Range("D22").Value = "'9.42"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("E23").Value = "  -2.38%  "

# Row 24
$ws.Range("D24").Value = "'77.34"
$ws.Range("E24").Value = "  -0.80%  "

# Row 25
$ws.Range("E25").Value = "  +8.47%  "

# Row 26
$ws.Range("D26").Value = "3.661.57"
$ws.Range("E26").Value = "  -0.52%  "

# Row 27
$ws.Range("E27").Value = "  +0.03%  "

# Row 28
$ws.Range("D28").Value = "'10.27"
$ws.Range("E28").Value = "  -2.13%  "

# Row 29
$ws.Range("E29").Value = "  -0.28%  "

# Row 30
$ws.Range("E30").Value = "  -2.48%  "

# Row 32
$ws.Range("E32").Value = "  -7.57%  "

# Row 33
$ws.Range("E33").Value = "  +4.38%  "

# Row 34
$ws.Range("E34").Value = "  -0.37%  "

# Row 35
$ws.Range("E35").Value = "  -1.00%  "

# Row 36
$ws.Range("D36").Value = "3.512.12"
$ws.Range("E36").Value = "  -0.72%  "

# Row 37
$ws.Range("E37").Value = "  -3.29%  "

# Row 38
$ws.Range("D38").Value = "'8.00"
$ws.Range("E38").Value = "  +0.35%  "

# Row 39
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'177.49"
$ws.Range("E40").Value = "  +0.69%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("E42").Value = "  +3.24%  "

# Row 43
$ws.Range("E43").Value = "  +0.28%  "

# Row 44
$ws.Range("D44").Value = "'5.43"
$ws.Range("E44").Value = "  -3.28%  "

# Row 45
$ws.Range("E45").Value = "  -0.60%  "

# Row 46
$ws.Range("D46").Value = "'28.47"
$ws.Range("E46").Value = "  -2.80%  "

# Row 47
$ws.Range("D47").Value = "'44.95"
$ws.Range("E47").Value = "  -1.93%  "

# Row 48
$ws.Range("E48").Value = "  +1.19%  "

# Row 49
$ws.Range("E49").Value = "  +3.46%  "

# Row 50
$ws.Range("E50").Value = "  -1.01%  "

# Row 51
$ws.Range("D51").Value = "'0.994"
$ws.Range("E51").Value = "  -1.65%  "
